$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 23
$ws.Range("B14").Value = "31/12/2025 02:46"
$ws.Range("C14").Value = 515
$ws.Range("D14").Value = "Conhecimentos Específicos"
$ws.Range("E14").Value = "Layout e Arranjos Físicos"
$ws.Range("F14").Value = "Rever tipos de layout"
